$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 1.33
$ws.Range("G5").Value = 2.37

$ws.Range("H27").Value = 3.7
$ws.Range("I27").Value = 2.1
$ws.Range("J27").Value = 3.5
$ws.Range("K27").Value = 2.4
$ws.Range("L27").Value = 2.63
$ws.Range("S27").Value = 1.29
$ws.Range("T27").Value = 3.5
$ws.Range("U27").Value = 1.5
$ws.Range("V27").Value = 2.5
$ws.Range("X27").Value = 19
$ws.Range("AB27").Value = 23
$ws.Range("AC27").Value = 17
$ws.Range("AG27").Value = 101
$ws.Range("AS27").Value = 101
$ws.Range("AT27").Value = 3.5
$ws.Range("AY27").Value = 17
$ws.Range("BA27").Value = 41

$ws.Range("G35").Value = 1.62
$ws.Range("H35").Value = 3.75
$ws.Range("I35").Value = 5.5
$ws.Range("J35").Value = 2.25
$ws.Range("U35").Value = 2
$ws.Range("V35").Value = 1.73
$ws.Range("Z35").Value = 12
$ws.Range("AB35").Value = 29
$ws.Range("AE35").Value = 19
$ws.Range("AK35").Value = 51
$ws.Range("AN35").Value = 3.5
$ws.Range("AO35").Value = 8.5
$ws.Range("AQ35").Value = 29
$ws.Range("AX35").Value = 29

$ws.Range("BD37").Value = 151

$ws.Range("R40").Value = 1.33

$ws.Range("Q41").Value = 1.72

$ws.Range("Q42").Value = 1.69
$ws.Range("R42").Value = 2.07
